# Refresh cryptocurrency price/volume data (cryptos list update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.908.57"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.630.32"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'206.57"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").Value = "'0.5133"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.2545"
$ws.Range("E8").Value = "  -3.50%  "
$ws.Range("D9").Value = "'0.06160"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").Value = "'20.28"
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("D11").Value = "'0.07542"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.630.62"
$ws.Range("E12").Value = "  -4.44%  "
$ws.Range("D13").Value = "'4.330"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "1.849.67"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "'0.5350"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "0.0₅7945"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "'64.93"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "25.931.75"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'4.592"
$ws.Range("E20").Value = "  -3.94%  "
$ws.Range("D21").Value = "'185.06"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'9.938"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "'6.040"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "'146.65"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'0.1190"
$ws.Range("E26").Value = "  -4.23%  "
$ws.Range("D27").Value = "'7.254"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("D28").Value = "'15.40"
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").Value = "'1.351"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "'0.05973"
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "'3.382"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'3.335"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").Value = "'0.9608"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("D36").Value = "'2.382"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").Value = "'0.5773"
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("D39").Value = "'0.01577"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").Value = "1.069.60"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").Value = "'5.763"
$ws.Range("D42").Value = "'1.002"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "'0.8393"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("D44").Value = "'99.55"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").Value = "1.783.30"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "'0.9985"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'53.92"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").Value = "'7.917"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'0.4225"
$ws.Range("E51").Value = "  -0.60%  "
